# error solve ifrs list
# Corrects the IFRS financial figures (columns D:AJ) for rows 2-9 of the
# "company_list" sheet, replacing the previously mis-scaled consolidated
# totals with the correct per-period values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 27585
$ws.Cells.Item(2, 5).Value = 7582
$ws.Cells.Item(2, 6).Value = 7582
$ws.Cells.Item(2, 7).Value = 6588
$ws.Cells.Item(2, 8).Value = 4518
$ws.Cells.Item(2, 9).Value = 4545
$ws.Cells.Item(2, 10).Value = -27
$ws.Cells.Item(2, 11).Value = 33945
$ws.Cells.Item(2, 12).Value = 15970
$ws.Cells.Item(2, 13).Value = 17974
$ws.Cells.Item(2, 14).Value = 17960
$ws.Cells.Item(2, 15).Value = 14
$ws.Cells.Item(2, 16).Value = 165
$ws.Cells.Item(2, 17).Value = 8366
$ws.Cells.Item(2, 18).Value = -4362
$ws.Cells.Item(2, 19).Value = -1119
$ws.Cells.Item(2, 20).Value = 2713
$ws.Cells.Item(2, 21).Value = 5653
$ws.Cells.Item(2, 22).Value = 4558
$ws.Cells.Item(2, 23).Value = 27.49
$ws.Cells.Item(2, 24).Value = 16.38
$ws.Cells.Item(2, 25).Value = 27.82
$ws.Cells.Item(2, 26).Value = 14.83
$ws.Cells.Item(2, 27).Value = 88.84999999999999
$ws.Cells.Item(2, 28).Value = 16483.45
$ws.Cells.Item(2, 29).Value = 2757
$ws.Cells.Item(2, 30).Value = 51.64
$ws.Cells.Item(2, 31).Value = 12190
$ws.Cells.Item(2, 32).Value = 11.68
$ws.Cells.Item(2, 33).Value = 156
$ws.Cells.Item(2, 34).Value = 0.11
$ws.Cells.Item(2, 35).Value = 5.07
$ws.Cells.Item(2, 36).Value = 164813395
# Row 3
$ws.Cells.Item(3, 4).Value = 32539
$ws.Cells.Item(3, 5).Value = 8302
$ws.Cells.Item(3, 6).Value = 7622
$ws.Cells.Item(3, 7).Value = 8232
$ws.Cells.Item(3, 8).Value = 5170
$ws.Cells.Item(3, 9).Value = 5187
$ws.Cells.Item(3, 10).Value = -18
$ws.Cells.Item(3, 11).Value = 43859
$ws.Cells.Item(3, 12).Value = 21163
$ws.Cells.Item(3, 13).Value = 22696
$ws.Cells.Item(3, 14).Value = 21245
$ws.Cells.Item(3, 15).Value = 1451
$ws.Cells.Item(3, 16).Value = 165
$ws.Cells.Item(3, 17).Value = 8514
$ws.Cells.Item(3, 18).Value = -7859
$ws.Cells.Item(3, 19).Value = 59
$ws.Cells.Item(3, 20).Value = 1208
$ws.Cells.Item(3, 21).Value = 7306
$ws.Cells.Item(3, 22).Value = 6677
$ws.Cells.Item(3, 23).Value = 25.52
$ws.Cells.Item(3, 24).Value = 15.89
$ws.Cells.Item(3, 25).Value = 26.46
$ws.Cells.Item(3, 26).Value = 13.29
$ws.Cells.Item(3, 27).Value = 93.25
$ws.Cells.Item(3, 28).Value = 19430.81
$ws.Cells.Item(3, 29).Value = 3147
$ws.Cells.Item(3, 30).Value = 41.81
$ws.Cells.Item(3, 31).Value = 14551
$ws.Cells.Item(3, 32).Value = 9.039999999999999
$ws.Cells.Item(3, 33).Value = 220
$ws.Cells.Item(3, 34).Value = 0.17
$ws.Cells.Item(3, 35).Value = 6.19
$ws.Cells.Item(3, 36).Value = 164813395
# Row 4
$ws.Cells.Item(4, 4).Value = 40226
$ws.Cells.Item(4, 5).Value = 11020
$ws.Cells.Item(4, 6).Value = 11020
$ws.Cells.Item(4, 7).Value = 11318
$ws.Cells.Item(4, 8).Value = 7591
$ws.Cells.Item(4, 9).Value = 7493
$ws.Cells.Item(4, 10).Value = 98
$ws.Cells.Item(4, 11).Value = 63706
$ws.Cells.Item(4, 12).Value = 22410
$ws.Cells.Item(4, 13).Value = 41296
$ws.Cells.Item(4, 14).Value = 35947
$ws.Cells.Item(4, 15).Value = 5348
$ws.Cells.Item(4, 16).Value = 165
$ws.Cells.Item(4, 17).Value = 11640
$ws.Cells.Item(4, 18).Value = -9417
$ws.Cells.Item(4, 19).Value = 6978
$ws.Cells.Item(4, 20).Value = 1538
$ws.Cells.Item(4, 21).Value = 10102
$ws.Cells.Item(4, 22).Value = 3773
$ws.Cells.Item(4, 23).Value = 27.4
$ws.Cells.Item(4, 24).Value = 18.87
$ws.Cells.Item(4, 25).Value = 26.2
$ws.Cells.Item(4, 26).Value = 14.11
$ws.Cells.Item(4, 27).Value = 54.27
$ws.Cells.Item(4, 28).Value = 30503.53
$ws.Cells.Item(4, 29).Value = 4546
$ws.Cells.Item(4, 30).Value = 34.09
$ws.Cells.Item(4, 31).Value = 24957
$ws.Cells.Item(4, 32).Value = 6.21
$ws.Cells.Item(4, 33).Value = 226
$ws.Cells.Item(4, 34).Value = 0.15
$ws.Cells.Item(4, 35).Value = 4.35
$ws.Cells.Item(4, 36).Value = 164813395
# Row 5
$ws.Cells.Item(5, 4).Value = 46785
$ws.Cells.Item(5, 5).Value = 11792
$ws.Cells.Item(5, 6).Value = 11792
$ws.Cells.Item(5, 7).Value = 11958
$ws.Cells.Item(5, 8).Value = 7701
$ws.Cells.Item(5, 9).Value = 7729
$ws.Cells.Item(5, 10).Value = -28
$ws.Cells.Item(5, 11).Value = 80193
$ws.Cells.Item(5, 12).Value = 27141
$ws.Cells.Item(5, 13).Value = 53052
$ws.Cells.Item(5, 14).Value = 47623
$ws.Cells.Item(5, 15).Value = 5428
$ws.Cells.Item(5, 16).Value = 165
$ws.Cells.Item(5, 17).Value = 9400
$ws.Cells.Item(5, 18).Value = -13103
$ws.Cells.Item(5, 19).Value = 6362
$ws.Cells.Item(5, 20).Value = 4722
$ws.Cells.Item(5, 21).Value = 4678
$ws.Cells.Item(5, 22).Value = 4156
$ws.Cells.Item(5, 23).Value = 25.2
$ws.Cells.Item(5, 24).Value = 16.46
$ws.Cells.Item(5, 25).Value = 18.5
$ws.Cells.Item(5, 26).Value = 10.7
$ws.Cells.Item(5, 27).Value = 51.16
$ws.Cells.Item(5, 28).Value = 36789.46
$ws.Cells.Item(5, 29).Value = 4689
$ws.Cells.Item(5, 30).Value = 37.1
$ws.Cells.Item(5, 31).Value = 32429
$ws.Cells.Item(5, 32).Value = 5.37
$ws.Cells.Item(5, 33).Value = 289
$ws.Cells.Item(5, 34).Value = 0.17
$ws.Cells.Item(5, 35).Value = 5.5
$ws.Cells.Item(5, 36).Value = 164813395
# Row 6
$ws.Cells.Item(6, 4).Value = 55869
$ws.Cells.Item(6, 5).Value = 9425
$ws.Cells.Item(6, 6).Value = 9425
$ws.Cells.Item(6, 7).Value = 11117
$ws.Cells.Item(6, 8).Value = 6279
$ws.Cells.Item(6, 9).Value = 6488
$ws.Cells.Item(6, 11).Value = 98812
$ws.Cells.Item(6, 12).Value = 39320
$ws.Cells.Item(6, 13).Value = 59491
$ws.Cells.Item(6, 14).Value = 52403
$ws.Cells.Item(6, 16).Value = 165
$ws.Cells.Item(6, 17).Value = 9735
$ws.Cells.Item(6, 18).Value = -3883
$ws.Cells.Item(6, 19).Value = 7510
$ws.Cells.Item(6, 20).Value = 5350
$ws.Cells.Item(6, 21).Value = 4385
$ws.Cells.Item(6, 22).Value = 11643
$ws.Cells.Item(6, 23).Value = 16.87
$ws.Cells.Item(6, 24).Value = 11.24
$ws.Cells.Item(6, 25).Value = 12.97
$ws.Cells.Item(6, 26).Value = 7.02
$ws.Cells.Item(6, 27).Value = 66.09
$ws.Cells.Item(6, 28).Value = 41077.9
$ws.Cells.Item(6, 29).Value = 3937
$ws.Cells.Item(6, 30).Value = 30.99
$ws.Cells.Item(6, 31).Value = 35847
$ws.Cells.Item(6, 32).Value = 3.4
$ws.Cells.Item(6, 33).Value = 314
$ws.Cells.Item(6, 34).Value = 0.26
$ws.Cells.Item(6, 35).Value = 7.07
$ws.Cells.Item(6, 36).Value = 164813395
# Row 7
$ws.Cells.Item(7, 4).Value = 65844
$ws.Cells.Item(7, 5).Value = 7630
$ws.Cells.Item(7, 7).Value = 7125
$ws.Cells.Item(7, 8).Value = 3530
$ws.Cells.Item(7, 9).Value = 4814
$ws.Cells.Item(7, 11).Value = 111885
$ws.Cells.Item(7, 12).Value = 47840
$ws.Cells.Item(7, 13).Value = 64048
$ws.Cells.Item(7, 14).Value = 57197
$ws.Cells.Item(7, 16).Value = 162
$ws.Cells.Item(7, 17).Value = 8530
$ws.Cells.Item(7, 18).Value = -9068
$ws.Cells.Item(7, 19).Value = 1053
$ws.Cells.Item(7, 20).Value = 4089
$ws.Cells.Item(7, 21).Value = 3944
$ws.Cells.Item(7, 23).Value = 11.59
$ws.Cells.Item(7, 24).Value = 5.36
$ws.Cells.Item(7, 25).Value = 8.789999999999999
$ws.Cells.Item(7, 26).Value = 3.35
$ws.Cells.Item(7, 27).Value = 74.69
$ws.Cells.Item(7, 29).Value = 2921
$ws.Cells.Item(7, 30).Value = 62.65
$ws.Cells.Item(7, 31).Value = 39325
$ws.Cells.Item(7, 32).Value = 4.65
$ws.Cells.Item(7, 33).Value = 287
$ws.Cells.Item(7, 34).Value = 0.16
$ws.Cells.Item(7, 35).Value = 9.82
# Row 8
$ws.Cells.Item(8, 4).Value = 76255
$ws.Cells.Item(8, 5).Value = 10492
$ws.Cells.Item(8, 7).Value = 10632
$ws.Cells.Item(8, 8).Value = 6284
$ws.Cells.Item(8, 9).Value = 7386
$ws.Cells.Item(8, 11).Value = 124675
$ws.Cells.Item(8, 12).Value = 53798
$ws.Cells.Item(8, 13).Value = 70877
$ws.Cells.Item(8, 14).Value = 64653
$ws.Cells.Item(8, 16).Value = 162
$ws.Cells.Item(8, 17).Value = 11201
$ws.Cells.Item(8, 18).Value = -6928
$ws.Cells.Item(8, 19).Value = 1280
$ws.Cells.Item(8, 20).Value = 3896
$ws.Cells.Item(8, 21).Value = 6273
$ws.Cells.Item(8, 23).Value = 13.76
$ws.Cells.Item(8, 24).Value = 8.24
$ws.Cells.Item(8, 25).Value = 12.12
$ws.Cells.Item(8, 26).Value = 5.32
$ws.Cells.Item(8, 27).Value = 75.90000000000001
$ws.Cells.Item(8, 29).Value = 4481
$ws.Cells.Item(8, 30).Value = 40.06
$ws.Cells.Item(8, 31).Value = 44451
$ws.Cells.Item(8, 32).Value = 4.04
$ws.Cells.Item(8, 33).Value = 336
$ws.Cells.Item(8, 34).Value = 0.19
$ws.Cells.Item(8, 35).Value = 7.5
# Row 9
$ws.Cells.Item(9, 4).Value = 85061
$ws.Cells.Item(9, 5).Value = 13533
$ws.Cells.Item(9, 7).Value = 14171
$ws.Cells.Item(9, 8).Value = 8732
$ws.Cells.Item(9, 9).Value = 9872
$ws.Cells.Item(9, 11).Value = 136904
$ws.Cells.Item(9, 12).Value = 56987
$ws.Cells.Item(9, 13).Value = 79929
$ws.Cells.Item(9, 14).Value = 74222
$ws.Cells.Item(9, 16).Value = 162
$ws.Cells.Item(9, 17).Value = 13016
$ws.Cells.Item(9, 18).Value = -6093
$ws.Cells.Item(9, 19).Value = 1421
$ws.Cells.Item(9, 20).Value = 3693
$ws.Cells.Item(9, 21).Value = 8231
$ws.Cells.Item(9, 23).Value = 15.91
$ws.Cells.Item(9, 24).Value = 10.27
$ws.Cells.Item(9, 25).Value = 14.22
$ws.Cells.Item(9, 26).Value = 6.68
$ws.Cells.Item(9, 27).Value = 71.3
$ws.Cells.Item(9, 29).Value = 5990
$ws.Cells.Item(9, 30).Value = 29.97
$ws.Cells.Item(9, 31).Value = 51030
$ws.Cells.Item(9, 32).Value = 3.52
$ws.Cells.Item(9, 33).Value = 369
$ws.Cells.Item(9, 34).Value = 0.21
$ws.Cells.Item(9, 35).Value = 6.16
